$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.235.12"
$ws.Range("E2").Value = "  +0.33%  "

# Row 3
$ws.Range("D3").Value = "1.604.90"
$ws.Range("E3").Value = "  +0.10%  "

# Row 4
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").Value = "212.77"
$ws.Range("E5").Value = "  +0.00%  "

# Row 6
$ws.Range("E6").Value = "  -0.08%  "

# Row 7
$ws.Range("D7").Value = "0.486"
$ws.Range("E7").Value = "  -0.21%  "

# Row 8
$ws.Range("E8").Value = "  +0.45%  "

# Row 9
$ws.Range("E9").Value = "  -0.43%  "

# Row 10
$ws.Range("D10").Value = "18.42"
$ws.Range("E10").Value = "  +2.15%  "

# Row 11
$ws.Range("D11").Value = "0.0812"
$ws.Range("E11").Value = "  -0.57%  "

# Row 12
$ws.Range("D12").Value = "1.827.84"
$ws.Range("E12").Value = "  +0.07%  "

# Row 13
$ws.Range("D13").Value = "1.600.62"
$ws.Range("E13").Value = "  -0.23%  "

# Row 14
$ws.Range("D14").Value = "4.02"
$ws.Range("E14").Value = "  +0.28%  "

# Row 15
$ws.Range("D15").Value = "0.513"
$ws.Range("E15").Value = "  +0.17%  "

# Row 16
$ws.Range("D16").Value = "26.209.03"

# Row 17
$ws.Range("D17").Value = "62.06"
$ws.Range("E17").Value = "  +2.58%  "

# Row 18
$ws.Range("D18").Value = "0.0₃0728"
$ws.Range("E18").Value = "  +0.69%  "

# Row 19
$ws.Range("E19").Value = "  -0.04%  "

# Row 20
$ws.Range("D20").Value = "200.58"
$ws.Range("E20").Value = "  -2.21%  "

# Row 21
$ws.Range("E21").Value = "  +0.54%  "

# Row 22
$ws.Range("E22").Value = "  -0.02%  "

# Row 23
$ws.Range("D23").Value = "6.00"
$ws.Range("E23").Value = "  +0.03%  "

# Row 24
$ws.Range("E24").Value = "  +2.54%  "

# Row 25
$ws.Range("D25").Value = "144.17"
$ws.Range("E25").Value = "  +1.64%  "

# Row 26
$ws.Range("E26").Value = "  -0.03%  "

# Row 27
$ws.Range("E27").Value = "  -2.62%  "

# Row 28
$ws.Range("D28").Value = "15.19"

# Row 29
$ws.Range("E29").Value = "  +1.93%  "

# Row 30
$ws.Range("E30").Value = "  +3.89%  "

# Row 31
$ws.Range("E31").Value = "  +0.46%  "

# Row 32
$ws.Range("E32").Value = "  +2.38%  "

# Row 33
$ws.Range("E33").Value = "  -1.40%  "

# Row 34
$ws.Range("E34").Value = "  +0.87%  "

# Row 35
$ws.Range("E35").Value = "  +1.25%  "

# Row 36
$ws.Range("D36").Value = "1.165.71"

# Row 37
$ws.Range("E37").Value = "  +2.56%  "

# Row 38
$ws.Range("E38").Value = "  -0.06%  "

# Row 39
$ws.Range("E39").Value = "  -0.55%  "

# Row 40
$ws.Range("E40").Value = "  +0.30%  "

# Row 42
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "5.34"
$ws.Range("E42").Value = "  +4.22%  "

# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "0.784"
$ws.Range("E43").Value = "  +0.32%  "

# Row 44
$ws.Range("D44").Value = "1.739.35"
$ws.Range("E44").Value = "  +0.03%  "

# Row 45
$ws.Range("D45").Value = "91.86"
$ws.Range("E45").Value = "  -1.22%  "

# Row 46
$ws.Range("D46").Value = "0.0⁦0107"
$ws.Range("E46").Value = "  +16.05%  "

# Row 47
$ws.Range("E47").Value = "  +1.04%  "

# Row 48
$ws.Range("D48").Value = "54.13"
$ws.Range("E48").Value = "  +1.23%  "

# Row 49
$ws.Range("E49").Value = "  +0.02%  "

# Row 50
$ws.Range("E50").Value = "  -0.46%  "

# Row 51
$ws.Range("E51").Value = "  -0.10%  "
